$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price/volume snapshot (Price in column D, Volume(1h) in
# column E) to the latest scrape. Row 10/11 also swap: "Polygon" and "OKB"
# traded ranking positions, so their Coin/Link/Price/Volume cells are
# rewritten with each other's (updated) data.
#
# Column D holds plain decimal-looking price text (e.g. "21.10", "0.5198").
# Excel's Range.Value setter auto-coerces such strings into floating point
# numbers, silently losing the exact text (trailing zeros / full precision --
# "21.10" would become 21.1, "0.5198" would become 0.51980000000000004). For
# those cells we briefly mark the cell as Text-formatted, assign the value,
# then clear the format again so the cell ends up with no explicit style --
# matching the source, which never set an explicit style on these cells.
# Cells whose new price text already can't parse as a number (it uses two
# '.' separators, e.g. "28.505.73") don't need this and are set directly, as
# are the column E percentages (they carry surrounding spaces/a '%' sign so
# Excel always keeps them as text).

$ws.Range("D2").Value = "28.505.73"
$ws.Range("E2").Value = "  -0.09%  "

$ws.Range("D3").Value = "1.824.75"
$ws.Range("E3").Value = "  -0.05%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.88"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.08%  "

$ws.Range("E6").Value = "  -0.01%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5198"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -2.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3888"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.26%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08371"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +8.34%  "

$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.91"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.28%  "

$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.114"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.06%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.430"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.16%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.11"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.36%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.003"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.04%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.523"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.48%  "

$ws.Range("D16").Value = "1.821.74"
$ws.Range("E16").Value = "  -0.25%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001133"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +4.77%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.08"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06601"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.09%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.78"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.51%  "

$ws.Range("E21").Value = "  +0.02%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.069"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.09%  "

$ws.Range("D23").Value = "28.541.69"
$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.43"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.65%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.276"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.68%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "21.10"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.28"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.68%  "

$ws.Range("D28").Value = "2.032.26"
$ws.Range("E28").Value = "  -0.17%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.402"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.42%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.52"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.27%  "

$ws.Range("E31").Value = "  -2.55%  "

$ws.Range("E32").Value = "  -3.16%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.726"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.24%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07419"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.91%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.663"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.15%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2226"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.09%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02364"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.62%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.222"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.15%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.808"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.83%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6327"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.99%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.40"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.53%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.192"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.399"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.26%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.62"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.89%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.785"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.75%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5966"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.86%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "126.69"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.99%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.993"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.07%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.205"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.16%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06984"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.47%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "74.42"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.17%  "
